$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AB1 = "obs", AC1 = "units"
$ws.Range("AB1").Value = "obs"

# Data rows: AB = "dsig/dpT" (establishes shared-string order before "units")
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 28).Value = "dsig/dpT"
}

$ws.Range("AC1").Value = "units"

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 29).Value = "pb"
}

$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("AC2:AC15").Select()
